# Modified presentation to give application examples of sorting algorithm.
#
# Slide 2 ("What is the algorithm?") content placeholder gets two new
# bullet points (indented one level) inserted between the existing
# "Problem it solves..." and "Goes from front to back..." bullets, and
# the illustrative picture on that slide is moved to a new position.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)

# --- Content Placeholder 2: insert the two new sub-bullets -----------------
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 2 is currently "Goes from front to back of list, sorting
# elements on the way" - insert the two new paragraphs right before it,
# immediately after "Problem it solves: sorting a list of numbers".
$lastPara = $tr.Paragraphs(2, 1)

$newPara1 = $lastPara.InsertBefore("One application is a program that tracks the height of all animals and the program needs to track animal information from shortest height to tallest height.`r")
$newPara1.IndentLevel = 2

# Re-fetch the "Goes from front..." paragraph, now pushed to index 3,
# and insert the second new bullet right before it as well.
$lastPara = $tr.Paragraphs(3, 1)
$newPara2 = $lastPara.InsertBefore("Another example is organizing names(if it’s modified to track alphabetical instead of numerical order).`r")
$newPara2.IndentLevel = 2

# --- Picture 6: reposition on the slide -------------------------------------
# Shape.Left/.Top are in points; OOXML stores EMU (914400 EMU = 1 inch =
# 72 points, i.e. 12700 EMU per point).
$pic = $s.Shapes.Item(3)
$pic.Left = 8729025 / 12700
$pic.Top  = 4020344 / 12700
